$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear old data rows (2-10) before rewriting, since the row layout/content changes entirely
$ws.Range("A2:F10").Clear()

# Row 1
$ws.Range("A1").Value = "BaseType"
$ws.Range("B1").Value = "Qty"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Health"
$ws.Range("E1").Value = "CombatDice"
$ws.Range("F1").Value = "Description"

# Row 2
$ws.Range("A2").Value = "zombie"
$ws.Range("B2").Value = 6
$ws.Range("C2").Value = "Small Zombie"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = ":R:"

# Row 3
$ws.Range("A3").Value = "zombie"
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = "Large Zombie"
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = ":R: :R:"

# Row 4
$ws.Range("A4").Value = "lychenthrope"
$ws.Range("B4").Value = 6
$ws.Range("C4").Value = "Small Lychenthrope"
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = ":W:"

# Row 5
$ws.Range("A5").Value = "lychenthrope"
$ws.Range("B5").Value = 5
$ws.Range("C5").Value = "Large Lychenthrope"
$ws.Range("D5").Value = 3
$ws.Range("E5").Value = ":W: :W:"

# Row 6
$ws.Range("A6").Value = "horror"
$ws.Range("B6").Value = 8
$ws.Range("C6").Value = "Horror"
$ws.Range("D6").Value = 4
$ws.Range("E6").Value = ":G: :G:"

# Row 7
$ws.Range("A7").Value = "event"
$ws.Range("B7").Value = 2
$ws.Range("C7").Value = "Discover Horrific Scene"
$ws.Range("F7").Value = "Gain 1 :corruption:"

# Row 8
$ws.Range("A8").Value = "event"
$ws.Range("B8").Value = 1
$ws.Range("C8").Value = "Healer's Care"
$ws.Range("F8").Value = "Gain 2 :heart:"

# Row 9
$ws.Range("A9").Value = "event"
$ws.Range("B9").Value = 1
$ws.Range("C9").Value = "Spring Trap"
$ws.Range("F9").Value = "Lose 3 :heart:"

# Row 10
$ws.Range("A10").Value = "event"
$ws.Range("B10").Value = 1
$ws.Range("C10").Value = "Tainted Blessing"
$ws.Range("F10").Value = "Gain 1 :corruption:`nAll players gain 3 :coin:"

# Row 11
$ws.Range("A11").Value = "event"
$ws.Range("B11").Value = 1
$ws.Range("C11").Value = "Eat Hearty Meal"
$ws.Range("F11").Value = "Gain 3 :heart:`nor`nAll players gain 1 :heart:"

# Row 12
$ws.Range("A12").Value = "event"
$ws.Range("B12").Value = 1
$ws.Range("C12").Value = "Learn a Lesson"
$ws.Range("F12").Value = "Each player draws the top card of the Purchase deck"

# Row 13
$ws.Range("A13").Value = "event"
$ws.Range("B13").Value = 1
$ws.Range("C13").Value = "Receive Exorcism"
$ws.Range("F13").Value = "Remove 2 Corruption cards from your deck and / or discard pile.  Then shuffle your deck."

# Row 14
$ws.Range("A14").Value = "event"
$ws.Range("B14").Value = 2
$ws.Range("C14").Value = "Murder of Crows"
$ws.Range("F14").Value = "Reveal the next Doom card"

# Row 15
$ws.Range("A15").Value = "event"
$ws.Range("B15").Value = 1
$ws.Range("C15").Value = "Gossip Among Villiagers"
$ws.Range("F15").Value = "Reveal any Disturbance in the City"

# Row 16
$ws.Range("A16").Value = "event"
$ws.Range("B16").Value = 1
$ws.Range("C16").Value = "Find Cache of Goods"
$ws.Range("F16").Value = "Gain 3 :coin:"

# Update selection to match target (C2)
$ws.Range("C2").Select()
